# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.802.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4737"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2935"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06524"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.61"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7424"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.887.74"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.251"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "286.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.786.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007544"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.134.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.336"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.293"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.234"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.928"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09813"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.343"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.316"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.188"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04918"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7009"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01903"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.76%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.40%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.330"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.09"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.018"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4308"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8395"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.94"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.613"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.045"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "915.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3979"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.64%  "

